# ADD THE SECURITY PLAN
#
# Inserts a new "SECURITY PLAN" section right after the paragraph that
# holds the embedded Word-document OLE object ("HARDWARE LIST"),
# leaving the pre-existing trailing empty paragraph in place (it ends
# up as the blank line right after the "How the hardware will be
# powered..." heading line), and appends the rest of the security
# write-up (with a couple of blank separator paragraphs) after it.
#
# NOTE: paragraph / range objects returned by this host do not "track"
# later document edits (they're point-in-time snapshots), so every step
# below re-fetches the paragraph it needs to touch by its *current*
# index instead of reusing an older object across a mutation.

$d = $word.ActiveDocument

$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-DocXmlPackage([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        "<w:document $wordNs><w:body>$bodyXml</w:body></w:document>" +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# A bare, run-less "<w:p/>" fragment (used for the blank separator lines).
$blankParaXml = New-DocXmlPackage('<w:p/>')

# Rewrites the paragraph currently at index $idx into a clean, run-less
# empty paragraph ("<w:p/>"), in place.
function Set-ParagraphBlank([int]$idx) {
    $p = $d.Paragraphs.Item($idx)
    $full = $d.Range($p.Range.Start, $p.Range.End)
    $full.InsertXML($blankParaXml) | Out-Null
}

# --- Locate the anchor: the pre-existing trailing empty paragraph that
# sits right before <w:sectPr> (the last paragraph in the document body,
# "w14:paraId=5C98F276" in the source file). ---
$anchorIdx = $d.Paragraphs.Count

# === Paragraphs inserted BEFORE the anchor ================================

# 1) "SECURITY PLAN " needs a <w:lastRenderedPageBreak/> ahead of its text
#    run, so build it from a raw WordprocessingML fragment instead of just
#    assigning .Range.Text.
$d.Paragraphs.Item($anchorIdx).Range.InsertParagraphBefore() | Out-Null
$securityIdx = $anchorIdx
$anchorIdx = $anchorIdx + 1
$securityXml = New-DocXmlPackage(
    '<w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">SECURITY PLAN </w:t></w:r></w:p>'
)
$d.Paragraphs.Item($securityIdx).Range.InsertXML($securityXml) | Out-Null

# 2) "How the hardware will be powered and connected to the internet"
$d.Paragraphs.Item($anchorIdx).Range.InsertParagraphBefore() | Out-Null
$howIdx = $anchorIdx
$anchorIdx = $anchorIdx + 1
$d.Paragraphs.Item($howIdx).Range.Text = "How the hardware will be powered and connected to the internet"

# At this point paragraph $anchorIdx is exactly where the diff's
# unchanged "<w:p/>" line is: right after the "How the hardware..."
# paragraph and still the original, untouched paragraph.

# === Paragraphs inserted AFTER the anchor ==================================

$afterIdx = $anchorIdx

# 3) Power paragraph
$d.Paragraphs.Item($afterIdx).Range.InsertParagraphAfter() | Out-Null
$afterIdx = $afterIdx + 1
$d.Paragraphs.Item($afterIdx).Range.Text = "Power:  PiSugar 1200 mAh boosts the Raspberry Pi Zero, charges the battery over USB-C and exposes battery telemetry (voltage/percentage). The Pi reads that battery status via the PiSugar daemon/CLI."

# 4) blank separator
$d.Paragraphs.Item($afterIdx).Range.InsertParagraphAfter() | Out-Null
$afterIdx = $afterIdx + 1
$d.Paragraphs.Item($afterIdx).Range.Text = "x"
Set-ParagraphBlank $afterIdx

# 5) Sensors paragraph
$d.Paragraphs.Item($afterIdx).Range.InsertParagraphAfter() | Out-Null
$afterIdx = $afterIdx + 1
$d.Paragraphs.Item($afterIdx).Range.Text = "Sensors: A vibration motor and a buzzer are each driven from 5v rails through 2N2222 NPN transistors, with 1 kΩ base resistors, 1N4148 diodes and a 100 µF bulk capacitor across 5v GND near the motor (to damp inrush). Status LED is current-limited by 220–330 Ω resistor."

# 6) blank separator
$d.Paragraphs.Item($afterIdx).Range.InsertParagraphAfter() | Out-Null
$afterIdx = $afterIdx + 1
$d.Paragraphs.Item($afterIdx).Range.Text = "x"
Set-ParagraphBlank $afterIdx

# 7) Network paragraph
$d.Paragraphs.Item($afterIdx).Range.InsertParagraphAfter() | Out-Null
$afterIdx = $afterIdx + 1
$d.Paragraphs.Item($afterIdx).Range.Text = "Network: The Pi connects to the internet using Wi-Fi* (client mode) to the user’s phone hotspot (WPA2). No Bluetooth is required for data—Wi-Fi handles everything."

# 8) Offline note paragraph
$d.Paragraphs.Item($afterIdx).Range.InsertParagraphAfter() | Out-Null
$afterIdx = $afterIdx + 1
$d.Paragraphs.Item($afterIdx).Range.Text = "  (If there’s no internet, the device works offline and syncs later.)"

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
